# edit.ps1 - applies the diff to the single-slide presentation.
#
# Summary of changes:
#   1. Add a new full-slide white background rectangle as the first shape
#      in the z-order (id=2, name "正方形/長方形 1").
#   2. Widen the "秘密鍵..." text box (id=35) and rewrite its third
#      paragraph's text.
#   3. Delete the "特製！なんぼかまし暗号方式！" text box (id=9).
#   4. Nudge the arrow connector (id=12) by a couple EMU and switch its
#      line dash style to "long dash".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. New background rectangle, sent behind every other shape.
# ---------------------------------------------------------------------
$bg = $s.Shapes.AddShape(1, 0, 0, (12192000 / 12700), (6858000 / 12700))
$bg.Name = "正方形/長方形 1"
$bg.Fill.Solid()
$bg.Fill.ForeColor.SchemeColor = "bg1"
$bg.Line.Visible = $false
$bg.TextFrame.VerticalAnchor = 3
$bg.TextFrame.TextRange.ParagraphFormat.Alignment = 2
# Drop the placeholder empty run the engine seeds new shapes with, so the
# paragraph serialises with no <a:r> at all (matching a never-typed-into
# shape) instead of an empty one.
$bg.TextFrame.TextRange.Text = "x"
$bg.TextFrame.TextRange.Characters(1, 1).Delete()
$bg.ZOrder(1)

# ---------------------------------------------------------------------
# 2. Shape 35: widen the box and replace the third paragraph's text.
# ---------------------------------------------------------------------
$box35 = Get-ShapeById $s 35
$box35.Width = 5262979 / 12700

$tr35 = $box35.TextFrame.TextRange
$p3 = $tr35.Paragraphs(3, 1)
# First overwrite with unrelated text so the engine can't line up a
# "common prefix" with the old run and split it into two runs; then set
# the real text so it lands in a single fresh run.
$p3.Text = "ZZZZZZZZZZZZZZZZZZZZZZZZZZZZZZZZ"
$p3b = $tr35.Paragraphs(3, 1)
$p3b.Text = "秘密鍵不正使用防止＝＞信頼できる公開鍵リスト"

# ---------------------------------------------------------------------
# 3. Delete shape 9 ("特製！なんぼかまし暗号方式！").
# ---------------------------------------------------------------------
$box9 = Get-ShapeById $s 9
if ($box9 -ne $null) {
    $box9.Delete()
}

# ---------------------------------------------------------------------
# 4. Connector 12: tiny offset tweak + dashed line.
# ---------------------------------------------------------------------
$conn12 = Get-ShapeById $s 12
$conn12.Left = 2647930 / 12700
$conn12.Top = 861603 / 12700
$conn12.Line.DashStyle = 7
